# Updates cryptos list values (Price / Volume(1h) columns, and a couple of
# coin name/link swaps) to match the latest scrape, per commit message:
# "Updated cryptos list on Sun Apr 30 17:21:16 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (leading apostrophe forces Excel to
# treat numeric-looking strings like "1.007" as text instead of a number,
# matching the original inline-string cell content) and reset the style back
# to Normal so no stray number-format style gets attached to the cell.
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.955.04"
Set-TextValue "E2" "  +1.57%  "
Set-TextValue "D3" "1.940.69"
Set-TextValue "E3" "  +1.08%  "
Set-TextValue "D4" "1.007"
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "335.23"
Set-TextValue "E5" "  +2.72%  "
Set-TextValue "D6" "1.006"
Set-TextValue "E6" "  -0.15%  "
Set-TextValue "D7" "0.4842"
Set-TextValue "E7" "  +0.24%  "
Set-TextValue "D8" "0.4135"
Set-TextValue "E8" "  +1.37%  "
Set-TextValue "D9" "0.08181"
Set-TextValue "E9" "  -0.62%  "
Set-TextValue "D10" "1.017"
Set-TextValue "E10" "  -0.67%  "
Set-TextValue "D11" "23.73"
Set-TextValue "E11" "  +0.69%  "
Set-TextValue "D12" "1.949.65"
Set-TextValue "E12" "  +1.80%  "
Set-TextValue "D13" "6.101"
Set-TextValue "E13" "  +0.79%  "
Set-TextValue "E14" "  +1.03%  "
Set-TextValue "D15" "91.29"
Set-TextValue "E15" "  -0.11%  "
Set-TextValue "D16" "0.06858"
Set-TextValue "E16" "  +0.73%  "
Set-TextValue "D17" "1.009"
Set-TextValue "E17" "  -0.04%  "
Set-TextValue "E18" "  -0.32%  "
Set-TextValue "D19" "17.83"
Set-TextValue "E19" "  +0.05%  "
Set-TextValue "E20" "  -0.01%  "
Set-TextValue "D21" "29.928.26"
Set-TextValue "D22" "5.649"
Set-TextValue "D23" "11.90"
Set-TextValue "E23" "  +1.07%  "
Set-TextValue "D24" "2.186"
Set-TextValue "E24" "  -0.44%  "
Set-TextValue "D25" "2.182.46"
Set-TextValue "E25" "  +1.43%  "
Set-TextValue "D26" "6.665"
Set-TextValue "E26" "  +0.70%  "
Set-TextValue "D27" "156.65"
Set-TextValue "E27" "  -0.02%  "
Set-TextValue "D28" "20.10"
Set-TextValue "E28" "  -0.14%  "
Set-TextValue "D29" "2.108"
Set-TextValue "E29" "  -0.59%  "
Set-TextValue "D30" "121.39"
Set-TextValue "E30" "  +0.72%  "
Set-TextValue "D31" "1.014"
Set-TextValue "E31" "  -1.14%  "
Set-TextValue "D32" "0.09645"
Set-TextValue "E32" "  +0.66%  "
Set-TextValue "D33" "5.601"
Set-TextValue "E33" "  +0.94%  "
Set-TextValue "D34" "1.424"
Set-TextValue "E34" "  +2.96%  "
Set-TextValue "D35" "3.550"
Set-TextValue "E35" "  -0.39%  "
Set-TextValue "D36" "0.06591"
Set-TextValue "E36" "  +7.32%  "
Set-TextValue "D37" "0.02288"
Set-TextValue "E37" "  +0.06%  "
Set-TextValue "D38" "1.217"
Set-TextValue "E38" "  +3.12%  "
Set-TextValue "D39" "0.5969"
Set-TextValue "E39" "  -0.29%  "
Set-TextValue "E40" "  -0.15%  "
Set-TextValue "E41" "  -1.03%  "
Set-TextValue "D42" "0.1854"
Set-TextValue "E42" "  -0.15%  "
Set-TextValue "D43" "2.502"
Set-TextValue "E43" "  +3.51%  "
Set-TextValue "B44" "WEMIXToken"
Set-TextValue "C44" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D44" "1.238"
Set-TextValue "E44" "  -3.31%  "
Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "12.36"
Set-TextValue "E45" "  -0.84%  "
Set-TextValue "D46" "0.07517"
Set-TextValue "E46" "  -1.14%  "
Set-TextValue "D47" "0.5582"
Set-TextValue "E47" "  -0.01%  "
Set-TextValue "D48" "1.989"
Set-TextValue "E48" "  +1.51%  "
Set-TextValue "D49" "117.81"
Set-TextValue "E49" "  -0.02%  "
Set-TextValue "D50" "72.80"
Set-TextValue "E50" "  -0.04%  "
Set-TextValue "D51" "2.418"
Set-TextValue "E51" "  -0.36%  "
